$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the old "reading list" / lab-technique reference rows that used to
# live below the main schedule table (F17, F20:F24, F27, F29:F34).
# ---------------------------------------------------------------------------
$ws.Range("F17").ClearContents()
$ws.Range("F20:F34").ClearContents()

# ---------------------------------------------------------------------------
# Re-number the "Barron and Boulpaep" chapter readings in column G (a new
# Ch. 32 reading was added at the front, and the old Ch. 38/Ch. 40 readings
# were replaced by a Ch. 37/Ch. 39 progression).
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = "Barron and Boulpaep Ch. 32"
$ws.Range("G4").Value = "Barron and Boulpaep Ch. 33"
$ws.Range("G5").Value = "Barron and Boulpaep Ch. 34"
$ws.Range("G6").Value = "Barron and Boulpaep Ch. 35"
$ws.Range("G7").Value = "Barron and Boulpaep Ch. 36"
$ws.Range("G8").Value = "Barron and Boulpaep Ch. 37"
$ws.Range("G9").Value = "Barron and Boulpaep Ch. 39"

# ---------------------------------------------------------------------------
# Add the new "Material" (D) / "Subject" (E) columns describing each week's
# R / Seurat curriculum content.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Intro to R & Setup"
$ws.Range("E3").Value = "Education Website - ROC 1"

$ws.Range("D4").Value = "Uploading files from GEO and other Seurat Functions"
$ws.Range("E4").Value = "ROC 2 + HW 1"

$ws.Range("D5").Value = "Rmarkdown & the whole kidney "
$ws.Range("E5").Value = "ROC 4 + HW 2"

$ws.Range("D6").Value = "Ddataset integration"
$ws.Range("E6").Value = "Week 3 Materials_for_Seurat_v5 + HW 3"

$ws.Range("D7").Value = "Psuedobulk"
$ws.Range("E7").Value = "AL-Pseudobulk-Materials + HW 4"

$ws.Range("D8").Value = "Data Visualization & Interpretation"
$ws.Range("E8").Value = "??"

$ws.Range("D9").Value = "Manuscript Writing"
$ws.Range("E9").Value = "??"

# ---------------------------------------------------------------------------
# Misc view bookkeeping to mirror the reverted workbook state.
# ---------------------------------------------------------------------------
$ws.Range("E19").Select()
